$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 0.1055014252150202
$ws.Range("B3").Value = 0.1223643052292021
$ws.Range("H3").Value = 0.2278657304442223
$ws.Range("B4").Value = 0.1543637739269478
$ws.Range("H4").Value = 0.2598651991419679
$ws.Range("B5").Value = 0.1804955279417598
$ws.Range("H5").Value = 0.28599695315678
$ws.Range("B6").Value = 0.1949553919506408
$ws.Range("C6").Value = 0.01001588205362942
$ws.Range("D6").Value = 19.42857643401977
$ws.Range("E6").Value = 0.0140652266181283
$ws.Range("F6").Value = 0.175252416895114
$ws.Range("G6").Value = 0.2146583670061658
$ws.Range("H6").Value = 0.300456817165661
$ws.Range("B7").Value = 0.02133977971314975
$ws.Range("C7").ClearContents()
$ws.Range("D7").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("F7").ClearContents()
$ws.Range("G7").ClearContents()
$ws.Range("H7").Value = 0.1268412049281699
$ws.Range("B8").Value = 0.01956737982069295
$ws.Range("C8").ClearContents()
$ws.Range("D8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("F8").ClearContents()
$ws.Range("G8").ClearContents()
$ws.Range("H8").Value = 0.1250688050357131
$ws.Range("B9").Value = 0.01967619388672123
$ws.Range("C9").Value = 0.002473789470562971
$ws.Range("D9").Value = 2.231242862100606
$ws.Range("E9").Value = 0.01607087454981219
$ws.Range("F9").Value = 0.01480674496207632
$ws.Range("G9").Value = 0.02454564281136651
$ws.Range("H9").Value = 0.1251776191017414
$ws.Range("B10").Value = 0.01819615792027296
$ws.Range("C10").Value = 0.001706920798387876
$ws.Range("D10").Value = 2.357010790265366
$ws.Range("E10").Value = 0.007233561167894183
$ws.Range("F10").Value = 0.01484935047587101
$ws.Range("G10").Value = 0.02154296536467539
$ws.Range("H10").Value = 0.1236975831352931
$ws.Range("B11").Value = 0.0303489923793352
$ws.Range("H11").Value = 0.1358504175943554
$ws.Range("B12").Value = 0.04891898940062141
$ws.Range("H12").Value = 0.1544204146156416
$ws.Range("B13").Value = 0.06510632449043699
$ws.Range("H13").Value = 0.1706077497054572
$ws.Range("B14").Value = 0.07010579993015824
$ws.Range("H14").Value = 0.1756072251451784
$ws.Range("B15").Value = 0.07370156985864657
$ws.Range("H15").Value = 0.1792029950736667
$ws.Range("B16").Value = 0.07686203111217532
$ws.Range("H16").Value = 0.1823634563271955
$ws.Range("B17").Value = 0.0814021779321821
$ws.Range("H17").Value = 0.1869036031472023
$ws.Range("B18").Value = -0.1055014252150202
$ws.Range("C18").Value = 0.01028323408717581
$ws.Range("D18").Value = -17.0575072969968
$ws.Range("E18").Value = 0.03729676889284375
$ws.Range("F18").Value = -0.1257324288707526
$ws.Range("G18").Value = -0.08527042155928778
$ws.Range("B19").Value = 0.08375116361469076
$ws.Range("H19").Value = 0.1892525888297109
$ws.Range("B20").Value = 0.08702112339154074
$ws.Range("C20").Value = 0.007828776380724812
$ws.Range("D20").Value = 1018685867.564496
$ws.Range("E20").Value = 0.05548274623810262
$ws.Range("F20").Value = 0.07163032926453368
$ws.Range("G20").Value = 0.1024119175185478
$ws.Range("H20").Value = 0.1925225486065609
$ws.Range("B21").Value = 0.09187251963339944
$ws.Range("H21").Value = 0.1973739448484196
$ws.Range("B22").Value = 0.09512302478531502
$ws.Range("C22").Value = 0.007559825112236007
$ws.Range("D22").Value = 1419454891999.033
$ws.Range("E22").Value = 0.04498861270255855
$ws.Range("F22").Value = 0.08026888315052809
$ws.Range("G22").Value = 0.109977166420102
$ws.Range("H22").Value = 0.2006244500003352
$ws.Range("B23").Value = 0.09897289051688896
$ws.Range("H23").Value = 0.2044743157319091
$ws.Range("B24").Value = 0.09870430536441838
$ws.Range("H24").Value = 0.2042057305794385
$ws.Range("B25").Value = 0.1043448419642999
$ws.Range("C25").Value = 0.007971645305570303
$ws.Range("D25").Value = 20.84991429154281
$ws.Range("E25").Value = 0.05534011883655374
$ws.Range("F25").Value = 0.08866328690278645
$ws.Range("G25").Value = 0.1200263970258132
$ws.Range("H25").Value = 0.2098462671793201
$ws.Range("B26").Value = 0.1043956610411623
$ws.Range("C26").Value = 0.008426823957779054
$ws.Range("D26").Value = 1032466629595.513
$ws.Range("E26").Value = 0.06334589898020578
$ws.Range("F26").Value = 0.08775956981154608
$ws.Range("G26").Value = 0.1210317522707784
$ws.Range("H26").Value = 0.2098970862561825
$ws.Range("B27").Value = 0.112019832054445
$ws.Range("C27").Value = 0.007584406861219519
$ws.Range("D27").Value = 22.12649488213892
$ws.Range("E27").Value = 0.05512387243952326
$ws.Range("F27").Value = 0.09712583766759547
$ws.Range("G27").Value = 0.1269138264412949
$ws.Range("H27").Value = 0.2175212572694651
$ws.Range("B28").Value = 0.1158879859826247
$ws.Range("C28").Value = 0.007907565227186959
$ws.Range("D28").Value = 22.15229369508844
$ws.Range("E28").Value = 0.0842547485804046
$ws.Range("F28").Value = 0.1003567383376262
$ws.Range("G28").Value = 0.131419233627623
$ws.Range("H28").Value = 0.2213894111976449
$ws.Range("B29").Value = 0.0212273710948751
$ws.Range("C29").Value = 0.002800350282614289
$ws.Range("D29").Value = 2.719612804864009
$ws.Range("E29").Value = 0.01777532955774431
$ws.Range("F29").Value = 0.01567169048919264
$ws.Range("G29").Value = 0.02678305170055838
$ws.Range("H29").Value = 0.1267287963098953
